$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (H1, I1) - match formatting of existing header row (bold, centered)
$ws.Range("H1").Value = "Threat_total"
$ws.Range("I1").Value = "Threat_prc"
$ws.Range("H1:I1").Font.Bold = $true
$ws.Range("H1:I1").HorizontalAlignment = -4108

# New data values for Threat_total (H) and Threat_prc (I) columns
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = 71

$ws.Range("H3").Value = 20
$ws.Range("I3").Value = 53

$ws.Range("H4").Value = 12
$ws.Range("I4").Value = 75

$ws.Range("H5").Value = 17
$ws.Range("I5").Value = 61

$ws.Range("H6").Value = 13
$ws.Range("I6").Value = 72
